$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "250.74"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.96"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.437"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05679"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.412"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.376"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8155"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9280"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1440"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07515"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03126"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09358"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.566"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001589"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04762"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005787"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006377"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004999"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001030"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.700"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.195"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0003044"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04033"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006764"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002708"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007566"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005801"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01010"
